$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 previously held "FAIL" (shared string index 22); change it to "PASS"
$ws.Range("D3").Value = "PASS"

# D4 previously held "PASS" (shared string index 21); clear it to a blank cell
$ws.Range("D4").ClearContents()

# Update the sheet's active selection to D3:D5 (active cell D3)
$ws.Range("D3:D5").Select()
